$d = $word.ActiveDocument

function Set-ParagraphText($para, [string]$newText) {
    $r = $para.Range
    $bodyRange = $d.Range($r.Start, $r.End - 1)
    $bodyRange.Text = $newText
}

Set-ParagraphText $d.Paragraphs(1) "⚡️🚀המאמר היומי של מייק 23.06.24:⚡️🚀"
Set-ParagraphText $d.Paragraphs(2) "TextGrad: Automatic “Differentiation” via Text"
Set-ParagraphText $d.Paragraphs(3) "אני קצת שיכור אחרי כמה שוטים ובירות באירוע המגניב של one-shot אבל התמדה בסקירות יומיות גברה על כך. הסקירה של היום מדברת גישה ש״מטילה״(project) את שיטת מורד גרדיאנט (gradient descent או פשוט GD) למקרה שהמשתנה שאנו מפטמים לפיו זה הפרומפט ולא משקלי המודל (שנותרות קבועים). כמו שאתם זוכרים GD הסטנדרטי מזיזים בכיוון הגרדיאנט השלילי של פונקציית לוס (מחסירים ממשקלי המודל את הגדיאנט מוכפל בקצב למידה). "
Set-ParagraphText $d.Paragraphs(4) "ב-GD הגרדיאנט מחושב בצורה ברורה (לפחות מתמטית) כי פונקציית לוס הינה גזירה ביחס למשקלי המודל. ד״א בראייה ממוחשבת ניתן לגזור את פונקציית לוס לפי הקלט (תמונה) מאותה הסיבה - לפעמים עושים זאת כדי לבנות תמונה הממזערת את הלוס עבור קטגוריה מסוימת."
Set-ParagraphText $d.Paragraphs(5) "אבל איך לגזור את המודל ביחס לטקסט? הכוונה כאן לא לגזור את פונקציית לוס לפי הייצוגים של טוקני הקלט (זה דווקא אפשרי כמו במקרה של תמונה ונקרא soft prompting). אך כאן מדובר ב״גזירה״ אשכרה לפי הטקסט עצמו. כמובן שמבחינה מתמטית זה די בעייתי כי טקסט הוא משתנה דיסקרטי. "
Set-ParagraphText $d.Paragraphs(6) "המאמר הופך את ומחליף גזירה מתמטית על ידי מה ״פידבק של שכבה ח לשכבה n-1״ וכאן לא מדובר בשכבות של מודלי שפה אלא בשכבות של כלים שונים המפעילים ומופעלים על ידי מודלי שפה (נגיד rag או כמה אג'נטים). אז בכל שלב אנו שואלים מודל שפה (אם הוא מופעל) איך היה ניתן לשפר את הפרומפט בשלב שלהם כדי לשפר את התוצאה ומעבירים את הפידבק לשכבה הקודמת. כמובן שהאגרגציה של פידבקים מתחילה ה-llm בשכבה האחרונה של המערכת ושיש לו סוג של פונקציית לוס בתור ״שערוך של איכות התשובה״. ומכאן מתחילה האגרגציה."
Set-ParagraphText $d.Paragraphs(7) "אז textgrad זה בגדול פרופגציה של פידבק טקסטואלי ופחות טקסט אבל עדיין המאמר חמוד כי מאפשר מערכות מורכבות מלא מעט כלים המערבים llms."

# Append the new final paragraph with the updated arXiv link
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = "https://arxiv.org/abs/2406.07496"

Write-Output ("Paragraph count: " + $d.Paragraphs.Count)
